$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.713.80"
$ws.Range("E2").Value = "  +2.02%  "
$ws.Range("D3").Value = "2.112.48"
$ws.Range("E3").Value = "  +10.20%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'331.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.34%  "
$ws.Range("D6").Value = "'0.9995"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.5231"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.23%  "
$ws.Range("D8").Value = "'0.4402"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.89%  "
$ws.Range("D9").Value = "'0.09024"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.99%  "
$ws.Range("D10").Value = "'46.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.94%  "
$ws.Range("D11").Value = "'1.178"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.26%  "
$ws.Range("D12").Value = "'25.27"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.63%  "
$ws.Range("D13").Value = "2.111.46"
$ws.Range("E13").Value = "  +10.80%  "
$ws.Range("D14").Value = "'6.797"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.07%  "
$ws.Range("D15").Value = "'7.750"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.02%  "
$ws.Range("D16").Value = "'97.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.65%  "
$ws.Range("D17").Value = "'0.00001141"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.98%  "
$ws.Range("D18").Value = "'1.001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").Value = "'0.06623"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.80%  "
$ws.Range("E20").Value = "  +4.00%  "
$ws.Range("D21").Value = "'6.418"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.91%  "
$ws.Range("D22").Value = "'0.9993"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "30.832.73"
$ws.Range("E23").Value = "  +2.37%  "
$ws.Range("D24").Value = "'12.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.31%  "
$ws.Range("D25").Value = "2.356.71"
$ws.Range("E25").Value = "  +10.79%  "
$ws.Range("D26").Value = "'2.259"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.15%  "
$ws.Range("D27").Value = "'23.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.32%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "'163.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "'2.538"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.93%  "
$ws.Range("D30").Value = "'134.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.94%  "
$ws.Range("D31").Value = "'1.190"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.66%  "
$ws.Range("D32").Value = "'0.1070"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.57%  "
$ws.Range("D33").Value = "'6.233"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.49%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.535"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +28.15%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'3.910"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.10%  "
$ws.Range("D36").Value = "'0.02589"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.25%  "
$ws.Range("D37").Value = "'5.612"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.56%  "
$ws.Range("D38").Value = "'0.06748"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.19%  "
$ws.Range("E39").Value = "  +12.28%  "
$ws.Range("D40").Value = "'9.533"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.29%  "
$ws.Range("E41").Value = "  +4.67%  "
$ws.Range("D42").Value = "'0.6803"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.75%  "
$ws.Range("D43").Value = "'1.253"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.39%  "
$ws.Range("D44").Value = "'14.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.21%  "
$ws.Range("D45").Value = "'0.9987"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").Value = "'0.6348"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.14%  "
$ws.Range("D47").Value = "'2.252"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.66%  "
$ws.Range("D48").Value = "'1.288"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.30%  "
$ws.Range("E49").Value = "  +0.75%  "
$ws.Range("D50").Value = "'124.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.74%  "
$ws.Range("D51").Value = "'83.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.16%  "
